# Update the "general" sheet: add two new rows (row 6 and 7) describing the
# prior distributions for fluxes and thermodynamic quantities, pushing the
# previously-existing rows (Number of exp. conditions, Number of model
# structures, Number of particles, Compute robust fluxes, Final tolerance)
# down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# Insert two blank rows above the old row 6 ("Number of exp. conditions...").
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()

# --- Row 6: Prior distribution for fluxes -------------------------------
$a6 = $ws.Range("A6")
$a6.Value = "Prior distribution for fluxes (uniform or normal)"
$a6.HorizontalAlignment = -4131
$a6.VerticalAlignment = -4160
$a6.Borders.LineStyle = 1
$a6.Borders.Weight = 2
$a6.Font.Name = "Calibri"
$a6.Font.Size = 11
$a6.Font.Bold = $true

$b6 = $ws.Range("B6")
$b6.Value = "normal"
$b6.HorizontalAlignment = -4108
$b6.VerticalAlignment = -4107
$b6.Borders.LineStyle = 1
$b6.Borders.Weight = 2
$b6.Font.Name = "Calibri"
$b6.Font.Size = 11
$b6.Font.Bold = $false

$ws.Rows.Item(6).RowHeight = 13.8

# --- Row 7: Prior distribution for thermodynamic quantities -------------
$a7 = $ws.Range("A7")
$a7.Value = "Prior distribution for thermodynamic quantities (uniform or normal)"
$a7.HorizontalAlignment = -4131
$a7.VerticalAlignment = -4160
$a7.Borders.LineStyle = 1
$a7.Borders.Weight = 2
$a7.Font.Name = "Calibri"
$a7.Font.Size = 11
$a7.Font.Bold = $true

$b7 = $ws.Range("B7")
$b7.Value = "normal"
$b7.HorizontalAlignment = -4108
$b7.VerticalAlignment = -4107
$b7.Borders.LineStyle = 1
$b7.Borders.Weight = 2
$b7.Font.Name = "Calibri"
$b7.Font.Size = 11
$b7.Font.Bold = $false

$ws.Rows.Item(7).RowHeight = 13.8

# Make "general" the active sheet/tab and select A6:B7, matching the
# post-edit view state (activeTab back to the general sheet).
$ws.Activate()
$ws.Range("A6:B7").Select()
